$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRESY")

# Row 8
$ws.Range("D8").Value = 894300
$ws.Range("E8").Value = 705300
$ws.Range("F8").Value = 358400
$ws.Range("G8").Value = 129700
$ws.Range("H8").Value = 105600
$ws.Range("I8").Value = 80900
$ws.Range("J8").Value = 65600

# Row 9
$ws.Range("D9").Value = 561500
$ws.Range("E9").Value = 494300
$ws.Range("F9").Value = 210400
$ws.Range("G9").Value = 109400
$ws.Range("H9").Value = 89800
$ws.Range("I9").Value = 51000
$ws.Range("J9").Value = 91400

# Row 10
$ws.Range("D10").Value = 332800
$ws.Range("E10").Value = 211000
$ws.Range("F10").Value = 148000
$ws.Range("G10").Value = 20200
$ws.Range("H10").Value = 15800
$ws.Range("I10").Value = 30000
$ws.Range("J10").Value = -25800

# Row 14
$ws.Range("D14").Value = -40200
$ws.Range("E14").Value = -4100
$ws.Range("F14").Value = -9700
$ws.Range("G14").Value = -29600
$ws.Range("H14").Value = -26000
$ws.Range("J14").Value = -1900

# Row 15
$ws.Range("D15").Value = 34400
$ws.Range("E15").Value = 124500
$ws.Range("F15").Value = 18400
$ws.Range("G15").Value = 300

# Row 17
$ws.Range("D17").Value = 207000
$ws.Range("E17").Value = 515900
$ws.Range("F17").Value = -92400
$ws.Range("G17").Value = 65600
$ws.Range("H17").Value = 78400
$ws.Range("I17").Value = 55600
$ws.Range("J17").Value = 50800

# Row 18
$ws.Range("D18").Value = 687400
$ws.Range("E18").Value = 189400
$ws.Range("F18").Value = 450800
$ws.Range("G18").Value = 64100
$ws.Range("H18").Value = 27200
$ws.Range("I18").Value = 25400
$ws.Range("J18").Value = 14800

# Row 20
$ws.Range("D20").Value = -372700
$ws.Range("E20").Value = 194500
$ws.Range("F20").Value = -16900
$ws.Range("G20").Value = -19700
$ws.Range("H20").Value = -52000
$ws.Range("I20").Value = -10200
$ws.Range("J20").Value = -3700

# Row 21
$ws.Range("D21").Value = 400600
$ws.Range("E21").Value = 462800
$ws.Range("F21").Value = 469100
$ws.Range("G21").Value = 50200
$ws.Range("H21").Value = -18100
$ws.Range("I21").Value = 21500
$ws.Range("J21").Value = 16400

# Row 22
$ws.Range("D22").Value = 187800
$ws.Range("E22").Value = 300100
$ws.Range("F22").Value = 124100
$ws.Range("G22").Value = 20000
$ws.Range("H22").Value = 16400
$ws.Range("I22").Value = 10900
$ws.Range("J22").Value = 9400

# Row 23
$ws.Range("D23").Value = 126900
$ws.Range("E23").Value = 83700
$ws.Range("F23").Value = 309800
$ws.Range("G23").Value = 24300
$ws.Range("H23").Value = -41200
$ws.Range("I23").Value = 4300
$ws.Range("J23").Value = 1700

# Row 24
$ws.Range("D24").Value = 5300
$ws.Range("E24").Value = 62200
$ws.Range("F24").Value = 132700
$ws.Range("G24").Value = 7000
$ws.Range("H24").Value = -8900
$ws.Range("I24").Value = 800

# Row 26
$ws.Range("D26").Value = 121600
$ws.Range("E26").Value = 21400
$ws.Range("F26").Value = 177100
$ws.Range("G26").Value = 17400
$ws.Range("H26").Value = -32300
$ws.Range("I26").Value = 3500
$ws.Range("J26").Value = 1300

# Row 27
$ws.Range("D27").Value = -162600
$ws.Range("E27").Value = -59200
$ws.Range("F27").Value = 91400
$ws.Range("G27").Value = 2600
$ws.Range("H27").Value = -24500
$ws.Range("I27").Value = -600
$ws.Range("J27").Value = -500

# Row 29
$ws.Range("D29").Value = 286300
$ws.Range("E29").Value = 93900
$ws.Range("F29").Value = 18700

# Row 32
$ws.Range("D32").Value = 372700
$ws.Range("E32").Value = -194500
$ws.Range("F32").Value = 16900
$ws.Range("G32").Value = 19700
$ws.Range("H32").Value = 52000
$ws.Range("I32").Value = 10200
$ws.Range("J32").Value = 3700

# Row 33
$ws.Range("D33").Value = 123700
$ws.Range("E33").Value = 34700
$ws.Range("F33").Value = 110200
$ws.Range("G33").Value = 2600
$ws.Range("H33").Value = -24500
$ws.Range("I33").Value = -600
$ws.Range("J33").Value = -500

# Row 35
$ws.Range("D35").Value = 123700
$ws.Range("E35").Value = 34700
$ws.Range("F35").Value = 110200
$ws.Range("G35").Value = 2600
$ws.Range("H35").Value = -24500
$ws.Range("I35").Value = -600
$ws.Range("J35").Value = -500

# Row 41
$ws.Range("D41").Value = 1130200
$ws.Range("E41").Value = 581800
$ws.Range("F41").Value = 646700
$ws.Range("G41").Value = 14500
$ws.Range("H41").Value = 23000
$ws.Range("I41").Value = 24000
$ws.Range("J41").Value = 7400

# Row 42
$ws.Range("D42").Value = 880200
$ws.Range("E42").Value = 325500
$ws.Range("F42").Value = 250700
$ws.Range("G42").Value = 11600
$ws.Range("H42").Value = 11400
$ws.Range("I42").Value = 15200
$ws.Range("J42").Value = 11000

# Row 43
$ws.Range("D43").Value = 515000
$ws.Range("E43").Value = 383600
$ws.Range("F43").Value = 312300
$ws.Range("G43").Value = 43800
$ws.Range("H43").Value = 30100
$ws.Range("I43").Value = 37900
$ws.Range("J43").Value = 16600

# Row 44
$ws.Range("D44").Value = 192400
$ws.Range("E44").Value = 157000
$ws.Range("F44").Value = 107700
$ws.Range("G44").Value = 14600
$ws.Range("H44").Value = 14700
$ws.Range("I44").Value = 10500
$ws.Range("J44").Value = 10300

# Row 45
$ws.Range("D45").Value = 301500
$ws.Range("E45").Value = 120300
$ws.Range("F45").Value = 43300
$ws.Range("G45").Value = 20400
$ws.Range("H45").Value = 35200
$ws.Range("I45").Value = 4600
$ws.Range("J45").Value = 6000

# Row 46
$ws.Range("D46").Value = 2997500
$ws.Range("E46").Value = 1568200
$ws.Range("F46").Value = 1037300
$ws.Range("G46").Value = 96600
$ws.Range("H46").Value = 114400
$ws.Range("I46").Value = 75400
$ws.Range("J46").Value = 40900

# Row 47
$ws.Range("D47").Value = 1323400
$ws.Range("E47").Value = 464100
$ws.Range("F47").Value = 581500
$ws.Range("G47").Value = 105600
$ws.Range("H47").Value = 75800
$ws.Range("I47").Value = 51200
$ws.Range("J47").Value = 59500

# Row 48
$ws.Range("D48").Value = 5988400
$ws.Range("E48").Value = 3028300
$ws.Range("F48").Value = 4265900
$ws.Range("G48").Value = 146100
$ws.Range("H48").Value = 144100
$ws.Range("I48").Value = 409100
$ws.Range("J48").Value = 244400

# Row 49
$ws.Range("D49").Value = 651800
$ws.Range("E49").Value = 285400
$ws.Range("F49").Value = 542000
$ws.Range("G49").Value = 4000
$ws.Range("H49").Value = 4000
$ws.Range("I49").Value = 17700

# Row 52
$ws.Range("D52").Value = 299100
$ws.Range("E52").Value = 192700
$ws.Range("F52").Value = 171300
$ws.Range("G52").Value = 22800
$ws.Range("H52").Value = 23800
$ws.Range("I52").Value = 8200
$ws.Range("J52").Value = 17600

# Row 54
$ws.Range("D54").Value = 10482400
$ws.Range("E54").Value = 5538800
$ws.Range("F54").Value = 4580000
$ws.Range("G54").Value = 364700
$ws.Range("H54").Value = 362100
$ws.Range("I54").Value = 284700
$ws.Range("J54").Value = 236900

# Row 57
$ws.Range("D57").Value = 526800
$ws.Range("E57").Value = 433000
$ws.Range("F57").Value = 380800
$ws.Range("G57").Value = 15600
$ws.Range("H57").Value = 13500
$ws.Range("I57").Value = 45700
$ws.Range("J57").Value = 35400

# Row 58
$ws.Range("D58").Value = 938200
$ws.Range("E58").Value = 534200
$ws.Range("F58").Value = 1077600
$ws.Range("G58").Value = 56800
$ws.Range("H58").Value = 60500
$ws.Range("I58").Value = 35000

# Row 59
$ws.Range("D59").Value = 207000
$ws.Range("E59").Value = 207100
$ws.Range("F59").Value = 126500
$ws.Range("G59").Value = 30200
$ws.Range("H59").Value = 36100
$ws.Range("I59").Value = 22300
$ws.Range("J59").Value = 12600

# Row 60
$ws.Range("D60").Value = 1672000
$ws.Range("E60").Value = 1174300
$ws.Range("F60").Value = 1046000
$ws.Range("G60").Value = 102600
$ws.Range("H60").Value = 110100
$ws.Range("I60").Value = 61100
$ws.Range("J60").Value = 46300

# Row 61
$ws.Range("D61").Value = 5481800
$ws.Range("E61").Value = 2569900
$ws.Range("F61").Value = 2152000
$ws.Range("G61").Value = 133800
$ws.Range("H61").Value = 121900
$ws.Range("I61").Value = 96100
$ws.Range("J61").Value = 63500

# Row 62
$ws.Range("D62").Value = 1007900
$ws.Range("E62").Value = 666600
$ws.Range("F62").Value = 528600
$ws.Range("G62").Value = 24700
$ws.Range("H62").Value = 28300
$ws.Range("I62").Value = 24400
$ws.Range("J62").Value = 21300

# Row 66
$ws.Range("D66").Value = 9782100
$ws.Range("E66").Value = 5162400
$ws.Range("F66").Value = 4266600
$ws.Range("G66").Value = 319800
$ws.Range("H66").Value = 317400
$ws.Range("I66").Value = 227600
$ws.Range("J66").Value = 178200

# Row 72
$ws.Range("D72").Value = 424400
$ws.Range("E72").Value = 347700
$ws.Range("F72").Value = 284900
$ws.Range("G72").Value = 6600
$ws.Range("H72").Value = 600
$ws.Range("I72").Value = 26300
$ws.Range("J72").Value = 25300

# Row 76
$ws.Range("D76").Value = 700200
$ws.Range("E76").Value = 376300
$ws.Range("F76").Value = 313400
$ws.Range("G76").Value = 44900
$ws.Range("H76").Value = 44600
$ws.Range("I76").Value = 57100
$ws.Range("J76").Value = 58700

# Row 81
$ws.Range("D81").Value = 123700
$ws.Range("E81").Value = 34700
$ws.Range("F81").Value = 110200
$ws.Range("G81").Value = 2600
$ws.Range("H81").Value = -24500
$ws.Range("I81").Value = -600
$ws.Range("J81").Value = -500

# Row 83
$ws.Range("D83").Value = 86300
$ws.Range("E83").Value = 79400
$ws.Range("F83").Value = 35400
$ws.Range("G83").Value = 5900
$ws.Range("H83").Value = 6800
$ws.Range("I83").Value = 6400
$ws.Range("J83").Value = 5300

# Row 89
$ws.Range("D89").Value = 308000
$ws.Range("E89").Value = 208600
$ws.Range("F89").Value = 96800
$ws.Range("G89").Value = 11300
$ws.Range("H89").Value = 20300
$ws.Range("I89").Value = 14900
$ws.Range("J89").Value = 15300

# Row 91
$ws.Range("D91").Value = -55900
$ws.Range("E91").Value = -53100
$ws.Range("F91").Value = -13600
$ws.Range("G91").Value = -5100
$ws.Range("H91").Value = -3200
$ws.Range("I91").Value = -8000
$ws.Range("J91").Value = -5400

# Row 94
$ws.Range("D94").Value = -274600
$ws.Range("E94").Value = -55400
$ws.Range("F94").Value = 198200
$ws.Range("G94").Value = 20000
$ws.Range("H94").Value = -20300
$ws.Range("I94").Value = -2100
$ws.Range("J94").Value = -8100

# Row 96
$ws.Range("D96").Value = -21600
$ws.Range("E96").Value = -49000
$ws.Range("F96").Value = -5500
$ws.Range("G96").Value = -800
$ws.Range("H96").Value = -5600
$ws.Range("I96").Value = -5500
$ws.Range("J96").Value = -3200

# Row 100
$ws.Range("D100").Value = -52700
$ws.Range("E100").Value = 43600
$ws.Range("F100").Value = -106600
$ws.Range("G100").Value = -40800
$ws.Range("H100").Value = -10200
$ws.Range("J100").Value = -11000

# Row 101
$ws.Range("D101").Value = 324100
$ws.Range("E101").Value = 61700
$ws.Range("F101").Value = 120400
$ws.Range("G101").Value = 900
$ws.Range("H101").Value = 9300
$ws.Range("I101").Value = 900
$ws.Range("J101").Value = -1300

# Row 102
$ws.Range("D102").Value = 304800
$ws.Range("E102").Value = 258500
$ws.Range("F102").Value = 308800
$ws.Range("G102").Value = -8500
$ws.Range("H102").Value = -1000
$ws.Range("I102").Value = 13200
$ws.Range("J102").Value = -5100
